$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("C1").Value = "Unidade"
$ws.Range("D1").Value = "Descrição"

# Row 2: Tinta
$ws.Range("A2").Value = "Tinta"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "litros"
$ws.Range("D2").Value = ""

# Row 3: Cimento
$ws.Range("A3").Value = "Cimento"
$ws.Range("B3").Value = 0.8
$ws.Range("C3").Value = "kg"
$ws.Range("D3").Value = ""

# Row 4: Azulejos
$ws.Range("A4").Value = "Azulejos"
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = "unidades"
$ws.Range("D4").Value = ""

# Row 5: Areia
$ws.Range("A5").Value = "Areia"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = "m³"
$ws.Range("D5").Value = ""

# Row 6: Tijolos
$ws.Range("A6").Value = "Tijolos"
$ws.Range("B6").Value = 480
$ws.Range("C6").Value = "unidades"
$ws.Range("D6").Value = ""

# Row 7: Cimento
$ws.Range("A7").Value = "Cimento"
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = "kg"
$ws.Range("D7").Value = ""

# Row 8: Areia
$ws.Range("A8").Value = "Areia"
$ws.Range("B8").Value = 20
$ws.Range("C8").Value = "kg"
$ws.Range("D8").Value = ""

# Row 9: Pedra
$ws.Range("A9").Value = "Pedra"
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = "kg"
$ws.Range("D9").Value = ""
